# Fill in the T-SQL column (column E) of the DataTypes worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataTypes")

# Map of row number -> T-SQL type value for column E
$values = @{
    2  = "BIT"
    3  = "VARBINARY(0)"
    4  = "CHAR(0)"
    5  = "TEXT"
    6  = "VARCHAR(0)"
    7  = "INT"
    8  = "BIGINT"
    9  = "INT"
    10 = "SMALLINT"
    11 = "FLOAT"
    12 = "FLOAT"
    13 = "DECIMAL(0, 0)"
    14 = "DATE"
    15 = "TIME"
    16 = "DATETIME"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 5).Value = $values[$row]
}

# Widen column E to fit the new, longer T-SQL type names (mirrors double-
# clicking the column-E border to auto-fit after typing the values), and
# update the active selection to match the authored workbook state.
$ws.Columns.Item(5).ColumnWidth = 14.43
$ws.Range("E5").Select() | Out-Null

$wb.Save()
